# Auto-generated Excel COM-interop edit script
# Applies numeric value updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# per the target diff (Pandaemonium_Profits.xlsx scheduled-runner update).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98 (ALC)
$ws.Range("H98").Value = 5126.4287
$ws.Range("I98").Value = 5317.5
$ws.Range("J98").Value = 3980
$ws.Range("K98").Value = 5317.5
$ws.Range("L98").Value = 3980
$ws.Range("M98").Value = -3819.5
$ws.Range("N98").Value = -6976

# Row 122 (ALC)
$ws.Range("H122").Value = 5126.4287
$ws.Range("I122").Value = 5317.5
$ws.Range("J122").Value = 3980
$ws.Range("K122").Value = 15952.5
$ws.Range("L122").Value = 11940
$ws.Range("M122").Value = -13502.5
$ws.Range("N122").Value = -16840

# Row 134 (ALC)
$ws.Range("H134").Value = 67055.45
$ws.Range("J134").Value = 67055.45
$ws.Range("L134").Value = 67055.45
$ws.Range("N134").Value = -77195.45

# Row 135 (ALC)
$ws.Range("H135").Value = 55556490
$ws.Range("I135").Value = 25000662
$ws.Range("J135").Value = 142858850
$ws.Range("K135").Value = 225005958
$ws.Range("L135").Value = 1285729650
$ws.Range("M135").Value = -225003423
$ws.Range("N135").Value = -1285734720

$ws = $wb.Worksheets.Item("ARM")
# Row 88 (ARM)
$ws.Range("H88").Value = 5395.647
$ws.Range("I88").Value = 11921.2
$ws.Range("J88").Value = 2676.6667
$ws.Range("K88").Value = 11921.2
$ws.Range("L88").Value = 2676.6667
$ws.Range("M88").Value = -11515.2
$ws.Range("N88").Value = -3488.6667

# Row 91 (ARM)
$ws.Range("H91").Value = 5395.647
$ws.Range("I91").Value = 11921.2
$ws.Range("J91").Value = 2676.6667
$ws.Range("K91").Value = 11921.2
$ws.Range("L91").Value = 2676.6667
$ws.Range("M91").Value = -10517.2
$ws.Range("N91").Value = -5484.6667

$ws = $wb.Worksheets.Item("BSM")
# Row 82 (BSM)
$ws.Range("H82").Value = 19376.79
$ws.Range("I82").Value = 4226.1816
$ws.Range("J82").Value = 40208.875
$ws.Range("K82").Value = 4226.1816
$ws.Range("L82").Value = 40208.875
$ws.Range("M82").Value = -3843.1816
$ws.Range("N82").Value = -40974.875

# Row 85 (BSM)
$ws.Range("H85").Value = 19376.79
$ws.Range("I85").Value = 4226.1816
$ws.Range("J85").Value = 40208.875
$ws.Range("K85").Value = 4226.1816
$ws.Range("L85").Value = 40208.875
$ws.Range("M85").Value = -2900.1816
$ws.Range("N85").Value = -42860.875

# Row 86 (BSM)
$ws.Range("H86").Value = 2085.1458
$ws.Range("I86").Value = 1904.878
$ws.Range("J86").Value = 3141
$ws.Range("K86").Value = 1904.878
$ws.Range("L86").Value = 3141
$ws.Range("M86").Value = -781.8779999999999
$ws.Range("N86").Value = -5387

# Row 89 (BSM)
$ws.Range("H89").Value = 2085.1458
$ws.Range("I89").Value = 1904.878
$ws.Range("J89").Value = 3141
$ws.Range("K89").Value = 9524.389999999999
$ws.Range("L89").Value = 15705
$ws.Range("M89").Value = -3908.389999999999
$ws.Range("N89").Value = -26937

$ws = $wb.Worksheets.Item("CRP")
# Row 2 (CRP)
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = $null
$ws.Range("N2").Value = 0

# Row 6 (CRP)
$ws.Range("H6").Value = 262514990
$ws.Range("J6").Value = 9999
$ws.Range("L6").Value = 9999
$ws.Range("N6").Value = -10225

# Row 7 (CRP)
$ws.Range("H7").Value = 150
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null

# Row 17 (CRP)
$ws.Range("H17").Value = 9908
$ws.Range("I17").Value = 9908
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 9908
$ws.Range("L17").Value = $null
$ws.Range("N17").Value = 0
$ws.Range("M17").Value = -9734

# Row 50 (CRP)
$ws.Range("H50").Value = 22899.4
$ws.Range("J50").Value = 22899.4
$ws.Range("L50").Value = 22899.4
$ws.Range("N50").Value = -24149.4

# Row 51 (CRP)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = $null
$ws.Range("N51").Value = 0

# Row 59 (CRP)
$ws.Range("H59").Value = 21487.309
$ws.Range("I59").Value = 5104
$ws.Range("J59").Value = 24466.092
$ws.Range("K59").Value = 5104
$ws.Range("L59").Value = 24466.092
$ws.Range("M59").Value = -3959
$ws.Range("N59").Value = -26756.092

# Row 60 (CRP)
$ws.Range("H60").Value = 93333
$ws.Range("I60").Value = 93333
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 93333
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = $null
$ws.Range("N60").Value = -92822

# Row 61 (CRP)
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = $null
$ws.Range("N61").Value = 0

# Row 68 (CRP)
$ws.Range("H68").Value = 20545
$ws.Range("J68").Value = 20545
$ws.Range("L68").Value = 20545
$ws.Range("N68").Value = -22043

# Row 71 (CRP)
$ws.Range("H71").Value = 20545
$ws.Range("J71").Value = 20545
$ws.Range("L71").Value = 61635
$ws.Range("N71").Value = -69123

# Row 74 (CRP)
$ws.Range("H74").Value = 40314
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 40314
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = $null
$ws.Range("M74").Value = 40314
$ws.Range("N74").Value = -42062

# Row 77 (CRP)
$ws.Range("H77").Value = 40314
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 40314
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = $null
$ws.Range("M77").Value = 120942
$ws.Range("N77").Value = -129678

# Row 94 (CRP)
$ws.Range("H94").Value = 881
$ws.Range("I94").Value = 762
$ws.Range("K94").Value = 762
$ws.Range("M94").Value = -311

# Row 99 (CRP)
$ws.Range("H99").Value = 1307.1428
$ws.Range("I99").Value = 1009.0909
$ws.Range("K99").Value = 1009.0909
$ws.Range("M99").Value = 488.9091

# Row 122 (CRP)
$ws.Range("H122").Value = 2701.147
$ws.Range("I122").Value = 2468.5833
$ws.Range("J122").Value = 3259.3
$ws.Range("K122").Value = 7405.749899999999
$ws.Range("L122").Value = 9777.900000000001
$ws.Range("M122").Value = -4955.749899999999
$ws.Range("N122").Value = -14677.9

# Row 126 (CRP)
$ws.Range("H126").Value = 1307.1428
$ws.Range("I126").Value = 1009.0909
$ws.Range("K126").Value = 3027.2727
$ws.Range("M126").Value = -557.2727

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 12828511
$ws.Range("I5").Value = 360.54544
$ws.Range("J5").Value = 83383336
$ws.Range("K5").Value = 1081.63632
$ws.Range("L5").Value = 250150008
$ws.Range("M5").Value = -969.6363200000001
$ws.Range("N5").Value = -250150232

# Row 12 (CUL)
$ws.Range("H12").Value = 66667010
$ws.Range("I12").Value = 166666880
$ws.Range("J12").Value = 427.8889
$ws.Range("K12").Value = 500000640
$ws.Range("L12").Value = 1283.6667
$ws.Range("M12").Value = -500000467
$ws.Range("N12").Value = -1629.6667

# Row 68 (CUL)
$ws.Range("H68").Value = 2368.1304
$ws.Range("I68").Value = 835.5
$ws.Range("J68").Value = 4040.0908
$ws.Range("K68").Value = 2506.5
$ws.Range("L68").Value = 12120.2724
$ws.Range("M68").Value = -1695.5
$ws.Range("N68").Value = -13742.2724

# Row 70 (CUL)
$ws.Range("H70").Value = 4479.643
$ws.Range("I70").Value = 2055.75
$ws.Range("K70").Value = 6167.25
$ws.Range("M70").Value = -5852.25

# Row 71 (CUL)
$ws.Range("H71").Value = 2368.1304
$ws.Range("I71").Value = 835.5
$ws.Range("J71").Value = 4040.0908
$ws.Range("K71").Value = 7519.5
$ws.Range("L71").Value = 36360.8172
$ws.Range("M71").Value = -3463.5
$ws.Range("N71").Value = -44472.8172

# Row 73 (CUL)
$ws.Range("H73").Value = 4479.643
$ws.Range("I73").Value = 2055.75
$ws.Range("K73").Value = 6167.25
$ws.Range("M73").Value = -5075.25

# Row 107 (CUL)
$ws.Range("H107").Value = 706.0862
$ws.Range("J107").Value = 1741.0588
$ws.Range("L107").Value = 5223.1764
$ws.Range("N107").Value = -9063.1764

# Row 113 (CUL)
$ws.Range("H113").Value = 515.65515
$ws.Range("I113").Value = 532.6896400000001
$ws.Range("J113").Value = 481.5862
$ws.Range("K113").Value = 1598.06892
$ws.Range("L113").Value = 1444.7586
$ws.Range("M113").Value = 571.9310799999998
$ws.Range("N113").Value = -5784.7586

# Row 122 (CUL)
$ws.Range("H122").Value = 796.7368
$ws.Range("I122").Value = 214
$ws.Range("J122").Value = 1444.2222
$ws.Range("K122").Value = 1926
$ws.Range("L122").Value = 12997.9998
$ws.Range("M122").Value = 524
$ws.Range("N122").Value = -17897.9998

# Row 123 (CUL)
$ws.Range("H123").Value = 1206
$ws.Range("I123").Value = 1206
$ws.Range("K123").Value = 3618
$ws.Range("M123").Value = -1168

# Row 135 (CUL)
$ws.Range("H135").Value = 12828511
$ws.Range("I135").Value = 360.54544
$ws.Range("J135").Value = 83383336
$ws.Range("K135").Value = 3244.90896
$ws.Range("L135").Value = 750450024
$ws.Range("M135").Value = -709.9089599999998
$ws.Range("N135").Value = -750455094

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (GSM)
$ws.Range("H102").Value = 3549.9583
$ws.Range("I102").Value = 3438.0833
$ws.Range("J102").Value = 3661.8333
$ws.Range("K102").Value = 3438.0833
$ws.Range("L102").Value = 3661.8333
$ws.Range("M102").Value = -1816.0833
$ws.Range("N102").Value = -6905.8333

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 1544.6666
$ws.Range("I22").Value = 366.66666
$ws.Range("J22").Value = 2133.6667
$ws.Range("K22").Value = 366.66666
$ws.Range("L22").Value = 2133.6667
$ws.Range("M22").Value = -71.66665999999998
$ws.Range("N22").Value = -2723.6667

# Row 23 (LTW)
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null

# Row 27 (LTW)
$ws.Range("H27").Value = 1544.6666
$ws.Range("I27").Value = 366.66666
$ws.Range("J27").Value = 2133.6667
$ws.Range("K27").Value = 366.66666
$ws.Range("L27").Value = 2133.6667
$ws.Range("M27").Value = -259.66666
$ws.Range("N27").Value = -2347.6667

# Row 40 (LTW)
$ws.Range("H40").Value = 4350.5
$ws.Range("I40").Value = 3500
$ws.Range("J40").Value = 6335
$ws.Range("K40").Value = 3500
$ws.Range("L40").Value = 6335
$ws.Range("M40").Value = -3364
$ws.Range("N40").Value = -6607

# Row 46 (LTW)
$ws.Range("H46").Value = 819.9167
$ws.Range("I46").Value = 593.9
$ws.Range("K46").Value = 593.9
$ws.Range("M46").Value = -405.9

# Row 55 (LTW)
$ws.Range("H55").Value = 1333966.6
$ws.Range("I55").Value = 4000000
$ws.Range("J55").Value = 950
$ws.Range("K55").Value = 4000000
$ws.Range("L55").Value = 950
$ws.Range("M55").Value = -3999827
$ws.Range("N55").Value = -1296

# Row 132 (LTW)
$ws.Range("H132").Value = 2566.262
$ws.Range("I132").Value = 2020.0667
$ws.Range("J132").Value = 3931.75
$ws.Range("K132").Value = 6060.2001
$ws.Range("L132").Value = 11795.25
$ws.Range("M132").Value = -3530.2001
$ws.Range("N132").Value = -16855.25

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 1900
$ws.Range("K122").Value = 5700
$ws.Range("M122").Value = -3250

# Row 132 (WVR)
$ws.Range("H132").Value = 2140.9707
$ws.Range("I132").Value = 2229.1
$ws.Range("J132").Value = 1480
$ws.Range("K132").Value = 6687.299999999999
$ws.Range("L132").Value = 4440
$ws.Range("M132").Value = -4157.299999999999
$ws.Range("N132").Value = -9500
